# Updated: ut 01. 06. 2021
# Apply changes to DailyStats sheet:
#  1. Clear AgTests (F) and AgPosit (G) columns for rows 393-422 (data retracted)
#  2. Correct several AgTests/AgPosit values for rows 425-452

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear F393:G422 entirely (removes the stored values for those cells)
$ws.Range("F393:G422").ClearContents()

# 2. Apply corrected values
$ws.Range("F425").Value = 138265
$ws.Range("F428").Value = 102299
$ws.Range("F429").Value = 178347
$ws.Range("F431").Value = 170844
$ws.Range("F432").Value = 123492
$ws.Range("G432").Value = 429
$ws.Range("G442").Value = 172
$ws.Range("F443").Value = 106620
$ws.Range("F444").Value = 103444
$ws.Range("F447").Value = 67096
$ws.Range("F449").Value = 59028
$ws.Range("F450").Value = 89809
$ws.Range("G450").Value = 168
$ws.Range("F451").Value = 84067
$ws.Range("G451").Value = 111
$ws.Range("F452").Value = 73422
$ws.Range("G452").Value = 123
